# Swap the step contents between TC3 and TC4 in the "Registrar Autorizações de
# Pagamento" test case sheet.
#
# Previously:
#   TC3 (row 23-28): step 2 = "Chefe Dado um registro selecionado..." /
#                              "SYSTEM Atualiza a lista de registros..."
#   TC4 (row 31-36): step 2 = "Chefe Clica para realizar a autorização de
#                              pagamento." / "SYSTEM Apresenta a tela de
#                              Registrar Autorizações de Pagamento"
#
# After the edit, the two steps are swapped: TC3's step becomes the old TC4
# content, and TC4's step becomes the old TC3 content.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldTC3Step = $ws.Range("B28").Value2
$oldTC3Result = $ws.Range("D28").Value2
$oldTC4Step = $ws.Range("B36").Value2
$oldTC4Result = $ws.Range("D36").Value2

$ws.Range("B28").Value2 = $oldTC4Step
$ws.Range("D28").Value2 = $oldTC4Result
$ws.Range("B36").Value2 = $oldTC3Step
$ws.Range("D36").Value2 = $oldTC3Result
